$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

# --- Refresh Price (D) and Volume(1h) (E) columns with latest market data ---
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("E3").Value = "  +5.11%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +5.20%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  +5.11%  "
$ws.Range("E16").Value = "  +4.31%  "
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  +13.05%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  +2.28%  "

# Price column cells need an explicit text round-trip so Excel
# does not reinterpret dotted-thousands price strings as numbers
Set-TextValue "D2" "69.169.35"
Set-TextValue "D3" "3.928.21"
Set-TextValue "D5" "604.91"
Set-TextValue "D6" "164.38"
Set-TextValue "D7" "3.926.75"
Set-TextValue "D11" "6.39"
Set-TextValue "D12" "0.463"
Set-TextValue "D13" "37.14"
Set-TextValue "D14" "0.0000246"
Set-TextValue "D15" "4.585.80"
Set-TextValue "D16" "3.898.44"
Set-TextValue "D17" "69.240.89"
Set-TextValue "D18" "7.51"
Set-TextValue "D20" "17.18"
Set-TextValue "D21" "11.25"
Set-TextValue "D22" "488.89"
Set-TextValue "D23" "0.725"
Set-TextValue "D25" "84.53"
Set-TextValue "D27" "12.16"
Set-TextValue "D28" "10.13"
Set-TextValue "D31" "4.080.17"
Set-TextValue "D32" "32.56"
Set-TextValue "D33" "7.88"
Set-TextValue "D34" "2.39"
Set-TextValue "D35" "3.872.56"
Set-TextValue "D36" "0.107"
Set-TextValue "D40" "1.00"
Set-TextValue "D42" "3.05"
Set-TextValue "D43" "441.74"
Set-TextValue "D44" "2.00"
Set-TextValue "D45" "48.49"
Set-TextValue "D49" "2.852.12"
Set-TextValue "D50" "141.54"

# --- Rows 46-48 rotated in the source ranking: USDe -> row46,
# EnergySwap -> row47, Cosmos -> row48 (each with refreshed price/volume) ---
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D46" "1.00"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "27.88"
$ws.Range("E47").Value = "  +19.84%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D48" "8.46"
$ws.Range("E48").Value = "  +0.28%  "
